$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row (Session 13) - all data collected
$ws.Range("A14").Value = 13
$ws.Range("C14").Value = 30

# Update the active selection to reflect the next empty row, like Excel would after data entry
$ws.Range("C15").Select()
